$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the unlabeled data row (old row 13, holding the teacher name that
# has moved up into row 10 "Objetivos:"). This shifts rows 14-25 up to 13-24.
$ws.Rows.Item(13).Delete()

# --- Content updates on the now-shifted grid ---

# Row 10 "Objetivos:" B/C now holds the responsible-teacher info instead of
# the long objectives paragraph.
$ws.Range("B10:C10").Value = "1112574 - Inês Conceição Roberto"

# Row 13 "Programa resumido:" B/C now holds "Semestral".
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 "Programa:" B/C now holds the activation date.
$ws.Range("B15:C15").Value = "01/01/2018"

# Row 18 "Método:" B/C now holds the responsible-teacher info.
$ws.Range("B18:C18").Value = "1112574 - Inês Conceição Roberto"

# Row 19 "Critério:" B/C now holds the evaluation-method text.
$ws.Range("B19:C19").Value = "Os alunos serão avaliados formalmente por duas provas teóricas. A ponderação das notas será de 50% para cada avaliação, ou seja: Média do período letivo normal = (P1 + P2 )/2"

# Row 20 "Norma de recuperação:" B/C now holds the approval-criterion text.
$ws.Range("B20:C20").Value = "Serão aprovados os alunos que obtiverem média igual ou maior que 5,0."

# Row 21 "Bibliografia:" B/C now holds the recovery-norm text.
$ws.Range("B21:C21").Value = "Aos alunos que não obtiverem média igual ou maior que 5,0, será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2 Serão aprovados os alunos que obtiverem média igual ou maior que 5,0"
